$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(9, 8).Value = 618521.7
$ws.Cells.Item(9, 9).Value = 1082393.6
$ws.Cells.Item(9, 10).Value = 25.777779
$ws.Cells.Item(9, 11).Value = 1082393.6
$ws.Cells.Item(9, 12).Value = 25.777779
$ws.Cells.Item(9, 13).Value = -1082224.6
$ws.Cells.Item(9, 14).Value = -363.777779
$ws.Cells.Item(15, 8).Value = 1250.0658
$ws.Cells.Item(15, 9).Value = 1250.0658
$ws.Cells.Item(15, 11).Value = 3750.1974
$ws.Cells.Item(15, 13).Value = -3581.1974
$ws.Cells.Item(40, 8).Value = 2930
$ws.Cells.Item(40, 9).Value = 2193.8
$ws.Cells.Item(40, 11).Value = 2193.8
$ws.Cells.Item(40, 13).Value = -2018.8
$ws.Cells.Item(43, 8).Value = 2919.6
$ws.Cells.Item(43, 9).Value = 2750
$ws.Cells.Item(43, 10).Value = 3032.6667
$ws.Cells.Item(43, 11).Value = 2750
$ws.Cells.Item(43, 12).Value = 3032.6667
$ws.Cells.Item(43, 13).Value = -2681
$ws.Cells.Item(43, 14).Value = -3170.6667
$ws.Cells.Item(51, 8).Value = 12214.286
$ws.Cells.Item(51, 10).Value = 12214.286
$ws.Cells.Item(51, 12).Value = 12214.286
$ws.Cells.Item(51, 14).Value = -13182.286
$ws.Cells.Item(58, 8).Value = 2629.5557
$ws.Cells.Item(58, 9).Value = 888.7
$ws.Cells.Item(58, 10).Value = 4805.625
$ws.Cells.Item(58, 11).Value = 2666.1
$ws.Cells.Item(58, 12).Value = 14416.875
$ws.Cells.Item(58, 13).Value = -2516.1
$ws.Cells.Item(58, 14).Value = -14716.875
$ws.Cells.Item(62, 8).Value = 5941.2
$ws.Cells.Item(62, 9).Value = 4500
$ws.Cells.Item(62, 10).Value = 6301.5
$ws.Cells.Item(62, 11).Value = 4500
$ws.Cells.Item(62, 12).Value = 6301.5
$ws.Cells.Item(62, 13).Value = -3876
$ws.Cells.Item(62, 14).Value = -7549.5
$ws.Cells.Item(64, 8).Value = 5561.2
$ws.Cells.Item(64, 9).Value = 6449.5
$ws.Cells.Item(64, 11).Value = 6449.5
$ws.Cells.Item(64, 13).Value = -6201.5
$ws.Cells.Item(65, 8).Value = 5941.2
$ws.Cells.Item(65, 9).Value = 4500
$ws.Cells.Item(65, 10).Value = 6301.5
$ws.Cells.Item(65, 11).Value = 22500
$ws.Cells.Item(65, 12).Value = 31507.5
$ws.Cells.Item(65, 13).Value = -19380
$ws.Cells.Item(65, 14).Value = -37747.5
$ws.Cells.Item(67, 8).Value = 5561.2
$ws.Cells.Item(67, 9).Value = 6449.5
$ws.Cells.Item(67, 11).Value = 6449.5
$ws.Cells.Item(67, 13).Value = -5591.5
$ws.Cells.Item(76, 8).Value = 5758.2
$ws.Cells.Item(76, 9).Value = 5947.75
$ws.Cells.Item(76, 11).Value = 5947.75
$ws.Cells.Item(76, 13).Value = -5632.75
$ws.Cells.Item(79, 8).Value = 5758.2
$ws.Cells.Item(79, 9).Value = 5947.75
$ws.Cells.Item(79, 11).Value = 5947.75
$ws.Cells.Item(79, 13).Value = -4855.75
$ws.Cells.Item(88, 8).Value = 3188.4375
$ws.Cells.Item(88, 10).Value = 5502.875
$ws.Cells.Item(88, 12).Value = 5502.875
$ws.Cells.Item(88, 14).Value = -6314.875
$ws.Cells.Item(91, 8).Value = 3188.4375
$ws.Cells.Item(91, 10).Value = 5502.875
$ws.Cells.Item(91, 12).Value = 5502.875
$ws.Cells.Item(91, 14).Value = -8310.875
$ws.Cells.Item(100, 8).Value = 5238.4
$ws.Cells.Item(100, 9).Value = 3830
$ws.Cells.Item(100, 11).Value = 3830
$ws.Cells.Item(100, 13).Value = -3289
$ws.Cells.Item(112, 8).Value = 7123.1333
$ws.Cells.Item(112, 10).Value = 7123.1333
$ws.Cells.Item(112, 12).Value = 21369.3999
$ws.Cells.Item(112, 14).Value = -23585.3999
$ws.Cells.Item(125, 8).Value = 424.7143
$ws.Cells.Item(125, 9).Value = 372.8
$ws.Cells.Item(125, 10).Value = 554.5
$ws.Cells.Item(125, 11).Value = 3355.2
$ws.Cells.Item(125, 12).Value = 4990.5
$ws.Cells.Item(125, 13).Value = -895.2000000000003
$ws.Cells.Item(125, 14).Value = -9910.5
$ws.Cells.Item(137, 8).Value = 4095.2222
$ws.Cells.Item(137, 9).Value = 4077.875
$ws.Cells.Item(137, 10).Value = 4109.1
$ws.Cells.Item(137, 11).Value = 12233.625
$ws.Cells.Item(137, 12).Value = 12327.3
$ws.Cells.Item(137, 13).Value = -9683.625
$ws.Cells.Item(137, 14).Value = -17427.3
$ws.Cells.Item(138, 8).Value = 5252.6714
$ws.Cells.Item(138, 9).Value = 5353.1665
$ws.Cells.Item(138, 10).Value = 5243.672
$ws.Cells.Item(138, 11).Value = 16059.4995
$ws.Cells.Item(138, 12).Value = 15731.016
$ws.Cells.Item(138, 13).Value = -10919.4995
$ws.Cells.Item(138, 14).Value = -26011.016
$ws.Cells.Item(141, 8).Value = 1937.1
$ws.Cells.Item(141, 9).Value = 1918.3572
$ws.Cells.Item(141, 11).Value = 5755.071599999999
$ws.Cells.Item(141, 13).Value = -575.0715999999993
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 2186.2144
$ws.Cells.Item(2, 9).Value = 2085.3635
$ws.Cells.Item(2, 11).Value = 2085.3635
$ws.Cells.Item(2, 13).Value = -1972.3635
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 13).Value = $null
$ws.Cells.Item(32, 8).Value = 18579.344
$ws.Cells.Item(32, 9).Value = 14813.905
$ws.Cells.Item(32, 10).Value = 36721.91
$ws.Cells.Item(32, 11).Value = 14813.905
$ws.Cells.Item(32, 12).Value = 36721.91
$ws.Cells.Item(32, 13).Value = -14526.905
$ws.Cells.Item(32, 14).Value = -37295.91
$ws.Cells.Item(45, 8).Value = 8974.333000000001
$ws.Cells.Item(45, 9).Value = 68666.664
$ws.Cells.Item(45, 10).Value = 3999.9722
$ws.Cells.Item(45, 11).Value = 68666.664
$ws.Cells.Item(45, 12).Value = 3999.9722
$ws.Cells.Item(45, 13).Value = -68289.664
$ws.Cells.Item(45, 14).Value = -4753.9722
$ws.Cells.Item(61, 8).Value = 5691.3125
$ws.Cells.Item(61, 9).Value = 3628
$ws.Cells.Item(61, 10).Value = 14632.333
$ws.Cells.Item(61, 11).Value = 3628
$ws.Cells.Item(61, 12).Value = 14632.333
$ws.Cells.Item(61, 13).Value = -3416
$ws.Cells.Item(61, 14).Value = -15056.333
$ws.Cells.Item(74, 8).Value = 5868.8
$ws.Cells.Item(74, 9).Value = 4385.5
$ws.Cells.Item(74, 10).Value = 9329.833000000001
$ws.Cells.Item(74, 11).Value = 4385.5
$ws.Cells.Item(74, 12).Value = 9329.833000000001
$ws.Cells.Item(74, 13).Value = -3511.5
$ws.Cells.Item(74, 14).Value = -11077.833
$ws.Cells.Item(76, 8).Value = 194368
$ws.Cells.Item(76, 10).Value = 194368
$ws.Cells.Item(76, 12).Value = 194368
$ws.Cells.Item(76, 14).Value = -195044
$ws.Cells.Item(77, 8).Value = 5868.8
$ws.Cells.Item(77, 9).Value = 4385.5
$ws.Cells.Item(77, 10).Value = 9329.833000000001
$ws.Cells.Item(77, 11).Value = 21927.5
$ws.Cells.Item(77, 12).Value = 46649.165
$ws.Cells.Item(77, 13).Value = -17559.5
$ws.Cells.Item(77, 14).Value = -55385.165
$ws.Cells.Item(79, 8).Value = 194368
$ws.Cells.Item(79, 10).Value = 194368
$ws.Cells.Item(79, 12).Value = 194368
$ws.Cells.Item(79, 14).Value = -196708
$ws.Cells.Item(110, 8).Value = 1728.4546
$ws.Cells.Item(110, 9).Value = 1768.4762
$ws.Cells.Item(110, 11).Value = 1768.4762
$ws.Cells.Item(110, 13).Value = 276.5237999999999
$ws.Cells.Item(116, 8).Value = 2186.2144
$ws.Cells.Item(116, 9).Value = 2085.3635
$ws.Cells.Item(116, 11).Value = 2085.3635
$ws.Cells.Item(116, 13).Value = 208.6365000000001
$ws.Cells.Item(122, 8).Value = 12802.444
$ws.Cells.Item(122, 9).Value = 11829.733
$ws.Cells.Item(122, 11).Value = 35489.199
$ws.Cells.Item(122, 13).Value = -33039.199
$ws.Cells.Item(132, 8).Value = 3096.4443
$ws.Cells.Item(132, 9).Value = 2985.5454
$ws.Cells.Item(132, 10).Value = 4316.3335
$ws.Cells.Item(132, 11).Value = 8956.636200000001
$ws.Cells.Item(132, 12).Value = 12949.0005
$ws.Cells.Item(132, 13).Value = -6426.636200000001
$ws.Cells.Item(132, 14).Value = -18009.0005
$ws.Cells.Item(136, 8).Value = 5691.3125
$ws.Cells.Item(136, 9).Value = 3628
$ws.Cells.Item(136, 10).Value = 14632.333
$ws.Cells.Item(136, 11).Value = 10884
$ws.Cells.Item(136, 12).Value = 43896.999
$ws.Cells.Item(136, 13).Value = -8334
$ws.Cells.Item(136, 14).Value = -48996.999
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 2186.2144
$ws.Cells.Item(3, 9).Value = 2085.3635
$ws.Cells.Item(3, 11).Value = 2085.3635
$ws.Cells.Item(3, 13).Value = -1971.3635
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).Value = $null
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).Value = $null
$ws.Cells.Item(86, 8).Value = 4030.7778
$ws.Cells.Item(86, 9).Value = 3781
$ws.Cells.Item(86, 10).Value = 4530.3335
$ws.Cells.Item(86, 11).Value = 3781
$ws.Cells.Item(86, 12).Value = 4530.3335
$ws.Cells.Item(86, 13).Value = -2658
$ws.Cells.Item(86, 14).Value = -6776.3335
$ws.Cells.Item(89, 8).Value = 4030.7778
$ws.Cells.Item(89, 9).Value = 3781
$ws.Cells.Item(89, 10).Value = 4530.3335
$ws.Cells.Item(89, 11).Value = 18905
$ws.Cells.Item(89, 12).Value = 22651.6675
$ws.Cells.Item(89, 13).Value = -13289
$ws.Cells.Item(89, 14).Value = -33883.6675
$ws.Cells.Item(105, 8).Value = 2271.2727
$ws.Cells.Item(105, 9).Value = 2553.7144
$ws.Cells.Item(105, 10).Value = 1777
$ws.Cells.Item(105, 11).Value = 2553.7144
$ws.Cells.Item(105, 12).Value = 1777
$ws.Cells.Item(105, 13).Value = -806.7143999999998
$ws.Cells.Item(105, 14).Value = -5271
$ws.Cells.Item(134, 8).Value = 2003.4546
$ws.Cells.Item(134, 9).Value = 1903.85
$ws.Cells.Item(134, 11).Value = 5711.549999999999
$ws.Cells.Item(134, 13).Value = -3176.549999999999
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 1545.2858
$ws.Cells.Item(7, 10).Value = 102.666664
$ws.Cells.Item(7, 12).Value = 102.666664
$ws.Cells.Item(7, 14).Value = -328.666664
$ws.Cells.Item(31, 8).Value = 6438.8184
$ws.Cells.Item(31, 9).Value = 5562.9624
$ws.Cells.Item(31, 10).Value = 8373
$ws.Cells.Item(31, 11).Value = 5562.9624
$ws.Cells.Item(31, 12).Value = 8373
$ws.Cells.Item(31, 13).Value = -5267.9624
$ws.Cells.Item(31, 14).Value = -8963
$ws.Cells.Item(34, 8).Value = 6438.8184
$ws.Cells.Item(34, 9).Value = 5562.9624
$ws.Cells.Item(34, 10).Value = 8373
$ws.Cells.Item(34, 11).Value = 5562.9624
$ws.Cells.Item(34, 12).Value = 8373
$ws.Cells.Item(34, 13).Value = -5360.9624
$ws.Cells.Item(34, 14).Value = -8777
$ws.Cells.Item(44, 8).Value = 28300
$ws.Cells.Item(44, 10).Value = 28300
$ws.Cells.Item(44, 12).Value = 28300
$ws.Cells.Item(44, 14).Value = -29184
$ws.Cells.Item(58, 8).Value = 7922
$ws.Cells.Item(58, 9).Value = 9082.666999999999
$ws.Cells.Item(58, 11).Value = 9082.666999999999
$ws.Cells.Item(58, 13).Value = -8879.666999999999
$ws.Cells.Item(62, 8).Value = 2886
$ws.Cells.Item(62, 9).Value = 2886
$ws.Cells.Item(62, 11).Value = 2886
$ws.Cells.Item(62, 13).Value = -2262
$ws.Cells.Item(65, 8).Value = 2886
$ws.Cells.Item(65, 9).Value = 2886
$ws.Cells.Item(65, 11).Value = 14430
$ws.Cells.Item(65, 13).Value = -11310
$ws.Cells.Item(86, 8).Value = 10707.833
$ws.Cells.Item(86, 9).Value = 13698.75
$ws.Cells.Item(86, 11).Value = 13698.75
$ws.Cells.Item(86, 13).Value = -12575.75
$ws.Cells.Item(88, 8).Value = 11062.25
$ws.Cells.Item(88, 10).Value = 10571.286
$ws.Cells.Item(88, 12).Value = 10571.286
$ws.Cells.Item(88, 14).Value = -11383.286
$ws.Cells.Item(89, 8).Value = 10707.833
$ws.Cells.Item(89, 9).Value = 13698.75
$ws.Cells.Item(89, 11).Value = 68493.75
$ws.Cells.Item(89, 13).Value = -62877.75
$ws.Cells.Item(91, 8).Value = 11062.25
$ws.Cells.Item(91, 10).Value = 10571.286
$ws.Cells.Item(91, 12).Value = 10571.286
$ws.Cells.Item(91, 14).Value = -13379.286
$ws.Cells.Item(105, 8).Value = 9328.808000000001
$ws.Cells.Item(105, 9).Value = 8597.869000000001
$ws.Cells.Item(105, 11).Value = 8597.869000000001
$ws.Cells.Item(105, 13).Value = -6850.869000000001
$ws.Cells.Item(112, 8).Value = 100000
$ws.Cells.Item(112, 10).Value = 100000
$ws.Cells.Item(112, 12).Value = 100000
$ws.Cells.Item(112, 14).Value = -102954
$ws.Cells.Item(122, 8).Value = 1400
$ws.Cells.Item(122, 9).Value = 1400
$ws.Cells.Item(122, 11).Value = 4200
$ws.Cells.Item(122, 13).Value = -1750
$ws.Cells.Item(132, 8).Value = 1929.2609
$ws.Cells.Item(132, 10).Value = 1096
$ws.Cells.Item(132, 12).Value = 3288
$ws.Cells.Item(132, 14).Value = -8348
$ws.Cells.Item(134, 8).Value = 2740
$ws.Cells.Item(134, 9).Value = 2423
$ws.Cells.Item(134, 10).Value = 3849.5
$ws.Cells.Item(134, 11).Value = 7269
$ws.Cells.Item(134, 12).Value = 11548.5
$ws.Cells.Item(134, 13).Value = -4734
$ws.Cells.Item(134, 14).Value = -16618.5
$ws.Cells.Item(136, 8).Value = 7922
$ws.Cells.Item(136, 9).Value = 9082.666999999999
$ws.Cells.Item(136, 11).Value = 27248.001
$ws.Cells.Item(136, 13).Value = -24698.001
$ws.Cells.Item(141, 8).Value = 173134.77
$ws.Cells.Item(141, 10).Value = 178060.16
$ws.Cells.Item(141, 12).Value = 178060.16
$ws.Cells.Item(141, 14).Value = -188420.16
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(7, 8).Value = 156.25
$ws.Cells.Item(7, 9).Value = 172.33333
$ws.Cells.Item(7, 10).Value = 140.16667
$ws.Cells.Item(7, 11).Value = 516.99999
$ws.Cells.Item(7, 12).Value = 420.50001
$ws.Cells.Item(7, 13).Value = -404.99999
$ws.Cells.Item(7, 14).Value = -644.50001
$ws.Cells.Item(12, 8).Value = 304.75
$ws.Cells.Item(12, 10).Value = 486.9091
$ws.Cells.Item(12, 12).Value = 1460.7273
$ws.Cells.Item(12, 14).Value = -1806.7273
$ws.Cells.Item(23, 8).Value = 308.4
$ws.Cells.Item(23, 10).Value = 308.4
$ws.Cells.Item(23, 12).Value = 925.1999999999999
$ws.Cells.Item(23, 14).Value = -1395.2
$ws.Cells.Item(39, 8).Value = 6749.75
$ws.Cells.Item(39, 9).Value = 1000
$ws.Cells.Item(39, 10).Value = 8666.333000000001
$ws.Cells.Item(39, 11).Value = 3000
$ws.Cells.Item(39, 12).Value = 25998.999
$ws.Cells.Item(39, 13).Value = -2706
$ws.Cells.Item(39, 14).Value = -26586.999
$ws.Cells.Item(46, 8).Value = 38725.777
$ws.Cells.Item(46, 9).Value = 1199.8695
$ws.Cells.Item(46, 11).Value = 3599.6085
$ws.Cells.Item(46, 13).Value = -3508.6085
$ws.Cells.Item(55, 8).Value = 7882.5557
$ws.Cells.Item(55, 10).Value = 10799.167
$ws.Cells.Item(55, 12).Value = 32397.501
$ws.Cells.Item(55, 14).Value = -32751.501
$ws.Cells.Item(68, 8).Value = 1811.3334
$ws.Cells.Item(68, 10).Value = 1811.3334
$ws.Cells.Item(68, 12).Value = 5434.0002
$ws.Cells.Item(68, 14).Value = -7056.0002
$ws.Cells.Item(71, 8).Value = 1811.3334
$ws.Cells.Item(71, 10).Value = 1811.3334
$ws.Cells.Item(71, 12).Value = 16302.0006
$ws.Cells.Item(71, 14).Value = -24414.0006
$ws.Cells.Item(117, 8).Value = 3324.5
$ws.Cells.Item(117, 9).Value = 2383
$ws.Cells.Item(117, 11).Value = 7149
$ws.Cells.Item(117, 13).Value = -3707
$ws.Cells.Item(131, 8).Value = 4579.4165
$ws.Cells.Item(131, 10).Value = 5195.5
$ws.Cells.Item(131, 12).Value = 15586.5
$ws.Cells.Item(131, 14).Value = -25666.5
$ws.Cells.Item(132, 8).Value = 2955.2222
$ws.Cells.Item(132, 9).Value = 2999.6667
$ws.Cells.Item(132, 10).Value = 2933
$ws.Cells.Item(132, 11).Value = 26997.0003
$ws.Cells.Item(132, 12).Value = 26397
$ws.Cells.Item(132, 13).Value = -24467.0003
$ws.Cells.Item(132, 14).Value = -31457
$ws.Cells.Item(138, 8).Value = 13184.375
$ws.Cells.Item(138, 9).Value = 3910
$ws.Cells.Item(138, 11).Value = 11730
$ws.Cells.Item(138, 13).Value = -6590
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 9).Value = 166.5
$ws.Cells.Item(2, 10).Value = 102.2
$ws.Cells.Item(2, 11).Value = 166.5
$ws.Cells.Item(2, 12).Value = 102.2
$ws.Cells.Item(2, 13).Value = -53.5
$ws.Cells.Item(2, 14).Value = -328.2
$ws.Cells.Item(18, 8).Value = 611703.4
$ws.Cells.Item(18, 9).Value = 611703.4
$ws.Cells.Item(18, 11).Value = 611703.4
$ws.Cells.Item(18, 13).Value = -611410.4
$ws.Cells.Item(46, 8).Value = 15791.615
$ws.Cells.Item(46, 10).Value = 15791.615
$ws.Cells.Item(46, 12).Value = 15791.615
$ws.Cells.Item(46, 14).Value = -16103.615
$ws.Cells.Item(52, 8).Value = 47499.5
$ws.Cells.Item(52, 9).Value = 47499.5
$ws.Cells.Item(52, 11).Value = 47499.5
$ws.Cells.Item(52, 13).Value = -47240.5
$ws.Cells.Item(57, 8).Value = 23755.445
$ws.Cells.Item(57, 10).Value = 36359
$ws.Cells.Item(57, 12).Value = 36359
$ws.Cells.Item(57, 14).Value = -37999
$ws.Cells.Item(80, 8).Value = 6686.2
$ws.Cells.Item(80, 10).Value = 8443.5
$ws.Cells.Item(80, 12).Value = 8443.5
$ws.Cells.Item(80, 14).Value = -10439.5
$ws.Cells.Item(83, 8).Value = 6686.2
$ws.Cells.Item(83, 10).Value = 8443.5
$ws.Cells.Item(83, 12).Value = 42217.5
$ws.Cells.Item(83, 14).Value = -52201.5
$ws.Cells.Item(105, 8).Value = 68250
$ws.Cells.Item(105, 10).Value = 68250
$ws.Cells.Item(105, 12).Value = 68250
$ws.Cells.Item(105, 14).Value = -75238
$ws.Cells.Item(107, 8).Value = 625
$ws.Cells.Item(107, 9).Value = 166.66667
$ws.Cells.Item(107, 11).Value = 166.66667
$ws.Cells.Item(107, 13).Value = 1753.33333
$ws.Cells.Item(122, 8).Value = 3999.8
$ws.Cells.Item(122, 9).Value = 3749.75
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 11249.25
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -8799.25
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(126, 8).Value = 4763.4165
$ws.Cells.Item(126, 9).Value = 4308.7144
$ws.Cells.Item(126, 10).Value = 5400
$ws.Cells.Item(126, 11).Value = 12926.1432
$ws.Cells.Item(126, 12).Value = 16200
$ws.Cells.Item(126, 13).Value = -10456.1432
$ws.Cells.Item(126, 14).Value = -21140
$ws.Cells.Item(132, 8).Value = 7862.213
$ws.Cells.Item(132, 9).Value = 7521.657
$ws.Cells.Item(132, 11).Value = 22564.971
$ws.Cells.Item(132, 13).Value = -20034.971
$ws.Cells.Item(135, 8).Value = 87745.31
$ws.Cells.Item(135, 10).Value = 87745.31
$ws.Cells.Item(135, 12).Value = 87745.31
$ws.Cells.Item(135, 14).Value = -97885.31
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 4214.231
$ws.Cells.Item(7, 9).Value = 3658.625
$ws.Cells.Item(7, 10).Value = 5103.2
$ws.Cells.Item(7, 11).Value = 3658.625
$ws.Cells.Item(7, 12).Value = 5103.2
$ws.Cells.Item(7, 13).Value = -3546.625
$ws.Cells.Item(7, 14).Value = -5327.2
$ws.Cells.Item(16, 8).Value = 1452.0714
$ws.Cells.Item(16, 9).Value = 1512.091
$ws.Cells.Item(16, 11).Value = 1512.091
$ws.Cells.Item(16, 13).Value = -1342.091
$ws.Cells.Item(22, 8).Value = 3177.6
$ws.Cells.Item(22, 9).Value = 2722
$ws.Cells.Item(22, 10).Value = 5000
$ws.Cells.Item(22, 11).Value = 2722
$ws.Cells.Item(22, 12).Value = 5000
$ws.Cells.Item(22, 13).Value = -2427
$ws.Cells.Item(22, 14).Value = -5590
$ws.Cells.Item(27, 8).Value = 3177.6
$ws.Cells.Item(27, 9).Value = 2722
$ws.Cells.Item(27, 10).Value = 5000
$ws.Cells.Item(27, 11).Value = 2722
$ws.Cells.Item(27, 12).Value = 5000
$ws.Cells.Item(27, 13).Value = -2615
$ws.Cells.Item(27, 14).Value = -5214
$ws.Cells.Item(32, 8).Value = 1706.5
$ws.Cells.Item(32, 9).Value = 1706.5
$ws.Cells.Item(32, 11).Value = 1706.5
$ws.Cells.Item(32, 13).Value = -1389.5
$ws.Cells.Item(40, 8).Value = 11030.4375
$ws.Cells.Item(40, 9).Value = 9781.182000000001
$ws.Cells.Item(40, 10).Value = 13778.8
$ws.Cells.Item(40, 11).Value = 9781.182000000001
$ws.Cells.Item(40, 12).Value = 13778.8
$ws.Cells.Item(40, 13).Value = -9645.182000000001
$ws.Cells.Item(40, 14).Value = -14050.8
$ws.Cells.Item(41, 8).Value = 20000
$ws.Cells.Item(41, 9).Value = 20000
$ws.Cells.Item(41, 11).Value = 20000
$ws.Cells.Item(41, 13).Value = -19562
$ws.Cells.Item(46, 8).Value = 10970.857
$ws.Cells.Item(46, 9).Value = 4665.6665
$ws.Cells.Item(46, 11).Value = 4665.6665
$ws.Cells.Item(46, 13).Value = -4477.6665
$ws.Cells.Item(55, 8).Value = 992.8333
$ws.Cells.Item(55, 9).Value = 988
$ws.Cells.Item(55, 10).Value = 993.8
$ws.Cells.Item(55, 11).Value = 988
$ws.Cells.Item(55, 12).Value = 993.8
$ws.Cells.Item(55, 13).Value = -815
$ws.Cells.Item(55, 14).Value = -1339.8
$ws.Cells.Item(82, 8).Value = 9159.799999999999
$ws.Cells.Item(82, 9).Value = 9949.5
$ws.Cells.Item(82, 10).Value = 8633.333000000001
$ws.Cells.Item(82, 11).Value = 9949.5
$ws.Cells.Item(82, 12).Value = 8633.333000000001
$ws.Cells.Item(82, 13).Value = -9588.5
$ws.Cells.Item(82, 14).Value = -9355.333000000001
$ws.Cells.Item(85, 8).Value = 9159.799999999999
$ws.Cells.Item(85, 9).Value = 9949.5
$ws.Cells.Item(85, 10).Value = 8633.333000000001
$ws.Cells.Item(85, 11).Value = 9949.5
$ws.Cells.Item(85, 12).Value = 8633.333000000001
$ws.Cells.Item(85, 13).Value = -8701.5
$ws.Cells.Item(85, 14).Value = -11129.333
$ws.Cells.Item(93, 8).Value = 904
$ws.Cells.Item(93, 9).Value = 866.2
$ws.Cells.Item(93, 10).Value = 998.5
$ws.Cells.Item(93, 11).Value = 866.2
$ws.Cells.Item(93, 12).Value = 998.5
$ws.Cells.Item(93, 13).Value = 381.8
$ws.Cells.Item(93, 14).Value = -3494.5
$ws.Cells.Item(106, 8).Value = 19703
$ws.Cells.Item(106, 10).Value = 19703
$ws.Cells.Item(106, 12).Value = 19703
$ws.Cells.Item(106, 14).Value = -22227
$ws.Cells.Item(109, 8).Value = 380061.66
$ws.Cells.Item(109, 10).Value = 380061.66
$ws.Cells.Item(109, 12).Value = 380061.66
$ws.Cells.Item(109, 14).Value = -382835.66
$ws.Cells.Item(122, 8).Value = 5354.091
$ws.Cells.Item(122, 9).Value = 5162.3335
$ws.Cells.Item(122, 10).Value = 5584.2
$ws.Cells.Item(122, 11).Value = 15487.0005
$ws.Cells.Item(122, 12).Value = 16752.6
$ws.Cells.Item(122, 13).Value = -13037.0005
$ws.Cells.Item(122, 14).Value = -21652.6
$ws.Cells.Item(126, 8).Value = 4214.231
$ws.Cells.Item(126, 9).Value = 3658.625
$ws.Cells.Item(126, 10).Value = 5103.2
$ws.Cells.Item(126, 11).Value = 10975.875
$ws.Cells.Item(126, 12).Value = 15309.6
$ws.Cells.Item(126, 13).Value = -8505.875
$ws.Cells.Item(126, 14).Value = -20249.6
$ws.Cells.Item(132, 8).Value = 2828.6843
$ws.Cells.Item(132, 9).Value = 2905.8572
$ws.Cells.Item(132, 10).Value = 2612.6
$ws.Cells.Item(132, 11).Value = 8717.571599999999
$ws.Cells.Item(132, 12).Value = 7837.799999999999
$ws.Cells.Item(132, 13).Value = -6187.571599999999
$ws.Cells.Item(132, 14).Value = -12897.8
$ws.Cells.Item(133, 8).Value = 90000
$ws.Cells.Item(133, 10).Value = 90000
$ws.Cells.Item(133, 12).Value = 90000
$ws.Cells.Item(133, 14).Value = -95060
$ws.Cells.Item(136, 8).Value = 4455.3477
$ws.Cells.Item(136, 9).Value = 3325.8
$ws.Cells.Item(136, 10).Value = 6573.25
$ws.Cells.Item(136, 11).Value = 9977.400000000001
$ws.Cells.Item(136, 12).Value = 19719.75
$ws.Cells.Item(136, 13).Value = -7427.400000000001
$ws.Cells.Item(136, 14).Value = -24819.75
$ws.Cells.Item(139, 8).Value = 87665.06
$ws.Cells.Item(139, 10).Value = 87665.06
$ws.Cells.Item(139, 12).Value = 87665.06
$ws.Cells.Item(139, 14).Value = -97945.06
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(80, 8).Value = 65554
$ws.Cells.Item(80, 10).Value = 65554
$ws.Cells.Item(80, 12).Value = 65554
$ws.Cells.Item(80, 14).Value = -67550
$ws.Cells.Item(81, 8).Value = 2660.3333
$ws.Cells.Item(81, 9).Value = 2744.125
$ws.Cells.Item(81, 11).Value = 5488.25
$ws.Cells.Item(81, 13).Value = -4427.25
$ws.Cells.Item(83, 8).Value = 65554
$ws.Cells.Item(83, 10).Value = 65554
$ws.Cells.Item(83, 12).Value = 196662
$ws.Cells.Item(83, 14).Value = -206646
$ws.Cells.Item(84, 8).Value = 2660.3333
$ws.Cells.Item(84, 9).Value = 2744.125
$ws.Cells.Item(84, 11).Value = 27441.25
$ws.Cells.Item(84, 13).Value = -22137.25
$ws.Cells.Item(96, 8).Value = 2899
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 14).Value = $null
$ws.Cells.Item(105, 8).Value = 55762.945
$ws.Cells.Item(105, 10).Value = 55762.945
$ws.Cells.Item(105, 12).Value = 55762.945
$ws.Cells.Item(105, 14).Value = -62750.945
$ws.Cells.Item(122, 8).Value = 5412.8667
$ws.Cells.Item(122, 9).Value = 4783.3076
$ws.Cells.Item(122, 10).Value = 9505
$ws.Cells.Item(122, 11).Value = 14349.9228
$ws.Cells.Item(122, 12).Value = 28515
$ws.Cells.Item(122, 13).Value = -11899.9228
$ws.Cells.Item(122, 14).Value = -33415
$ws.Cells.Item(132, 8).Value = 1736.2778
$ws.Cells.Item(132, 9).Value = 1661
$ws.Cells.Item(132, 11).Value = 4983
$ws.Cells.Item(132, 13).Value = -2453
